# Automatische test-sync: 2025-06-19 21:19:50
#
# Adds the new "Is product X op voorraad?" mail-log entry (row 10) to the
# "Logs" worksheet, extends the conditional-formatting ranges to cover the
# new row, and refreshes the "Dashboard" summary table so that the
# "Productinformatie" category (now counted twice) is resorted above
# "Openingstijden / Locatie" and "Factuur / Administratie".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new mail entry as row 10
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Is product X op voorraad?"
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Range("D10").Value = "Productinformatie"
$logs.Range("F10").Value = "2025-06-19 21:19:10"
$logs.Range("G10").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges from row 9 to row 10
# ---------------------------------------------------------------------
$dRules = $logs.Range("D2:D9").FormatConditions
for ($i = 1; $i -le $dRules.Count(); $i++) {
    $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D10"))
}

$gRules = $logs.Range("G2:G9").FormatConditions
for ($i = 1; $i -le $gRules.Count(); $i++) {
    $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G10"))
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: resort categories now that "Productinformatie"
#    counts 2 entries instead of 1
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 2

$dash.Range("A4").Value = "Openingstijden / Locatie"
$dash.Range("B4").Value = 1

$dash.Range("A5").Value = "Factuur / Administratie"
$dash.Range("B5").Value = 1
